$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.405.88"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "1.906.03"
$ws.Range("E3").Value = "  -2.90%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'239.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("D6").Value = "'0.9991"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.4726"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.68%  "
$ws.Range("D8").Value = "'0.2834"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.01%  "
$ws.Range("D9").Value = "'0.06656"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.81%  "
$ws.Range("D10").Value = "'18.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.06%  "
$ws.Range("D11").Value = "'100.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.72%  "
$ws.Range("D12").Value = "'0.07716"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").Value = "1.905.61"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").Value = "'5.197"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.15%  "
$ws.Range("D15").Value = "'0.6691"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "30.374.48"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").Value = "'254.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -9.03%  "
$ws.Range("D18").Value = "'0.9995"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").Value = "'0.000007443"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.75%  "
$ws.Range("D20").Value = "'12.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.87%  "
$ws.Range("D21").Value = "'5.369"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "'0.4489"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -12.86%  "
$ws.Range("D24").Value = "'6.309"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.50%  "
$ws.Range("D25").Value = "'9.352"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.56%  "
$ws.Range("D26").Value = "'165.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("D27").Value = "'18.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.71%  "
$ws.Range("D28").Value = "'2.047"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.48%  "
$ws.Range("D29").Value = "'0.1009"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.18%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.376"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.638"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Value = "'1.510"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.66%  "
$ws.Range("D33").Value = "'4.255"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.16%  "
$ws.Range("D34").Value = "'0.04711"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("D35").Value = "'0.7272"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("D36").Value = "'1.109"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.66%  "
$ws.Range("D37").Value = "'0.9985"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "'2.695"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("D39").Value = "'0.01914"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.95%  "
$ws.Range("E40").Value = "  -3.55%  "
$ws.Range("D41").Value = "'6.242"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.60%  "
$ws.Range("D42").Value = "'73.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.35%  "
$ws.Range("E43").Value = "  -8.93%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8593"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.17%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'105.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.81%  "
$ws.Range("D46").Value = "'0.9985"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "'0.4223"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.48%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.411"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.76%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'980.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("D50").Value = "'0.1194"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.73%  "
$ws.Range("D51").Value = "'34.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.49%  "
